# Thai translation pass for "Email 3 [TEMPLATE] Partner email - list of
# travel documents.docx"
#
# wdReplaceOne = 1 ; wdReplaceAll = 2 ; wdFindContinue = 1 ; wdCollapseEnd = 0

$d = $word.ActiveDocument

function Replace-All($find, $repl) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
        $true, 1, $false, $repl, 2) | Out-Null
}

# ---------------------------------------------------------------------
# 1. Language-switcher line at the very top of the document.
# ---------------------------------------------------------------------
Replace-All "English" "ภาษาอังกฤษ"
Replace-All " / Portuguese / French / Thai / Vietnamese / Spanish" " / ภาษาโปรตุเกส / ภาษาฝรั่งเศส /ภาษาไทย / ภาษาเวียดนาม / ภาษาสเปน"

# ---------------------------------------------------------------------
# 2. Brief table.
# ---------------------------------------------------------------------
Replace-All "Brief" "บทย่อ"

$findBrief = "An email sent to partners in the target country who have " + "RSVPed yes. We want them to submit their documents. It " + "will be sent via customer.io"
$replBrief = "อีเมล์ที่ส่งถึงคู่ค้าที่อยู่ในประเทศเป้าหมายและได้ตอบรับคำเชิญหรือ " + "RSVP แล้ว เราต้องการให้พวกเขาส่งเอกสารของพวกเขา โดยมันจะถูกส่งผ่านทาง customer.io"
Replace-All $findBrief $replBrief

Replace-All "Target audience" "กลุ่มเป้าหมาย"
Replace-All "Invited partners who RSVP yes" "พาร์ทเนอร์ที่ได้รับเชิญและได้รับตอบรับหรือ RSVP แล้ว"

# ---------------------------------------------------------------------
# 3. "Subject: [EVENT NAME] -- take the next step" line.
#    "Subject: " is split into a bold "หัวข้อ:" run and a separate plain
#    run containing just the space, to match the authored structure.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Subject: ") | Out-Null
$boldPart = $r.Duplicate
$boldPart.End = $boldPart.Start + 8
$spacePart = $r.Duplicate
$spacePart.Start = $spacePart.Start + 8
$spacePart.End = $r.End
$spacePart.Font.Bold = $false
$spacePart.Text = " "
$boldPart.Text = "หัวข้อ:"

Replace-All " — take the next step" " — ดำเนินการขั้นตอนต่อไป"

# ---------------------------------------------------------------------
# 4. Greeting block.
# ---------------------------------------------------------------------
Replace-All "Thank you for registering for " "ขอบคุณที่ลงทะเบียนสำหรับงาน "
Replace-All "Hi " "สวัสดี "

# Drop the "," run that used to follow [PARTNER NAME].
$r = $d.Content
$r.Find.Execute("[PARTNER NAME]") | Out-Null
$r.Collapse(0)
$r.MoveEnd(1, 1) | Out-Null
$r.Text = ""

Replace-All "We are excited for you to join us at " "เรารู้สึกตื่นเต้นที่คุณจะมาร่วมกิจกรรมกับเราที่งาน "

# ---------------------------------------------------------------------
# 5. Document checklist.
# ---------------------------------------------------------------------
$findConfirm = "To confirm your registration, we would require you and " + "one guest of your choice to provide us with:"
Replace-All $findConfirm "เพื่อยืนยันการลงทะเบียนของคุณ เราจะขอให้คุณและแขกผู้ติดตามหนึ่งคนที่คุณเลือกมานั้นช่วยส่งสิ่งต่อไปนี้ให้เรา:"

Replace-All "A signed copy of the " "สำเนาจรรยาบรรณและข้อกำหนดและเงื่อนไขที่ลงนามแล้ว "

# Remove the "Code of Conduct " run entirely (its comment anchor stays).
$r = $d.Content
$r.Find.Execute("Code of Conduct ") | Out-Null
$r.Text = ""

Replace-All "and " "และ "
Replace-All "Terms and Conditions" "ข้อกำหนดและเงื่อนไข"
Replace-All " (1 set from each person)" " (คนละ 1 ชุด)"
Replace-All "A scanned copy of your international passports" "สำเนาสแกนหนังสือเดินทางระหว่างประเทศของคุณ"
Replace-All "Covid-19 vaccination certificates" "ใบรับรองการฉีดวัคซีน Covid-19"

# ---------------------------------------------------------------------
# 6. Button + follow-up copy.
# ---------------------------------------------------------------------
Replace-All "Send my details" "ส่งรายละเอียดของฉัน"

$findManager = "Your country manager will be in touch to confirm your " + "booking or request any other relevant details. "
Replace-All $findManager "ผู้จัดการประจำประเทศของคุณจะติดต่อกับคุณเพื่อยืนยันการจองที่นั่งของคุณหรือเพื่อขอรายละเอียดอื่นๆ ที่เกี่ยวข้อง "

Replace-All "Our event package offers you and your guest: " "แพ็คเกจงานกิจกรรมของเราจะนำเสนอให้คุณและแขกผู้ติดตามของคุณดังนี้: "

# ---------------------------------------------------------------------
# 7. Event package bullet list.
# ---------------------------------------------------------------------
Replace-All "Flight tickets " "ตั๋วเครื่องบิน "
Replace-All "Travel insurance " "ประกันการเดินทาง "
Replace-All "Airport – Hotel – Airport transfer " "บริการรับ-ส่ง ระหว่าง สนามบิน – โรงแรม – สนามบิน "

$findRooms = "One hotel room for you and your guest / Two hotel rooms " + "for you and your guest"
Replace-All $findRooms "ห้องพักโรงแรมหนึ่งห้องสำหรับคุณและแขกของคุณ / ห้องพักโรงแรมสองห้องสำหรับคุณและแขกของคุณ"

Replace-All "Check-in" "เช็คอิน"
Replace-All "Check-out" "เช็คเอาท์"
Replace-All " on " " ในวันที่ "

Replace-All "Meals (Breakfast, lunch, and dinner)" "อาหาร (อาหารมื้อเช้า มื้อกลางวัน และมื้อเย็น)"
Replace-All "Sightseeing tour of " "ทัวร์เที่ยวชมเมือง "

$findLetter = "We will send you a confirmation letter before your " + "departure date with the event agenda and information " + "about your flights, transportation, and accommodation. "
Replace-All $findLetter "เราจะส่งจดหมายยืนยันให้คุณก่อนวันออกเดินทางพร้อมด้วยวาระการประชุมและข้อมูลเกี่ยวกับเที่ยวบิน การเดินทาง และที่พักของคุณ "

# ---------------------------------------------------------------------
# 8. Contact / sign-off block.
# ---------------------------------------------------------------------
Replace-All "If you have any questions, please contact us via " "หากคุณมีคำถามใดๆ กรุณาติดต่อเราผ่านทาง "
Replace-All "live chat" "แชทสด"

# " or " / " or " and ". " are ambiguous document-wide, so scope the Find
# to right after the anchor text that precedes each occurrence.
$r = $d.Content
$r.Find.Execute("แชทสด") | Out-Null
$r.Collapse(0)
$r.Find.Execute(" or ", $true, $false, $false, $false, $false, $true, 1, $false, " หรือทาง ", 1) | Out-Null

$r = $d.Content
$r.Find.Execute("WhatsApp") | Out-Null
$r.Collapse(0)
$r.Find.Execute(". ", $true, $false, $false, $false, $false, $true, 1, $false, " ", 1) | Out-Null

Replace-All "If you have any questions, please contact your country manager, " "หากคุณมีคำถามใดๆ โปรดติดต่อผู้จัดการประจำประเทศของคุณซึ่งได้แก่ "
Replace-All ", at " " ที่ "

$r = $d.Content
$r.Find.Execute("[EMAIL ADDRESS]") | Out-Null
$r.Collapse(0)
$r.Find.Execute(" or ", $true, $false, $false, $false, $false, $true, 1, $false, " หรือ ", 1) | Out-Null

Replace-All " (WhatsApp). " " (WhatsApp) "
Replace-All "We look forward to seeing you soon." "เราหวังว่าจะได้พบเจอคุณเร็วๆ นี้"

# ---------------------------------------------------------------------
# 9. Comments.
# ---------------------------------------------------------------------
foreach ($c in $d.Comments) {
    $t = $c.Range.Text
    if ($t -eq "link to T&C") {
        $c.Range.Text = "ลิงก์ไปยังข้อกำหนดและเงื่อนไข"
    } elseif ($t -eq "link to COC") {
        $c.Range.Text = "ลิงก์ไปยังจรรยาบรรณ"
    } elseif ($t -eq "please confirm these") {
        $c.Range.Text = "กรุณายืนยันสิ่งเหล่านี้"
    } elseif ($t -eq "choose either one") {
        $c.Range.Text = "เลือกอย่างใดอย่างหนึ่ง"
    } elseif ($t -eq "please check if these are all the required documents") {
        $c.Range.Text = "กรุณาตรวจสอบว่านี่คือเอกสารที่ต้องการทั้งหมดหรือไม่"
    }
}
